$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "extr1".."extr8" rows (8-15) down by two rows (to 10-17)
# so two fresh rows open up at 8-9 for the new "line7"/"line8" entries. Values
# are copied cell-by-cell (bottom row first) so we never overwrite data we
# still need to read.
for ($r = 15; $r -ge 8; $r--) {
    $dest = $r + 2
    $ws.Cells.Item($dest, 1).Value = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($dest, 2).Value = $ws.Cells.Item($r, 2).Value()
    $ws.Cells.Item($dest, 3).Value = $ws.Cells.Item($r, 3).Value()
    $ws.Cells.Item($dest, 4).Value = $ws.Cells.Item($r, 4).Value()
    $ws.Cells.Item($dest, 5).Value = $ws.Cells.Item($r, 5).Value()
}

# Rows 16/17 are brand new cells (never existed before), so column A needs the
# bold/bordered/centered style used by every other index cell. Pull that
# format from an existing styled cell (A2) rather than building it by hand,
# so we reuse the existing style record instead of minting a new one.
$ws.Cells.Item(2, 1).Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 8: line7 (new) ---
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

# --- Row 9: line8 (new) ---
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# --- Row 10: extr1 (was row 8) ---
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 3).Value = 5
$ws.Cells.Item(10, 4).Value = 12
$ws.Cells.Item(10, 5).Value = $true

# --- Row 11: extr2 (was row 9) ---
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 3).Value = 5
$ws.Cells.Item(11, 4).Value = 9
$ws.Cells.Item(11, 5).Value = $true

# --- Row 12: extr3 (was row 10) ---
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 3).Value = 10
$ws.Cells.Item(12, 4).Value = 11
$ws.Cells.Item(12, 5).Value = $true

# --- Row 13: extr4 (was row 11) ---
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 3).Value = 7
$ws.Cells.Item(13, 4).Value = 8
$ws.Cells.Item(13, 5).Value = $false

# --- Row 14: extr5 (was row 12) ---
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 11
$ws.Cells.Item(14, 5).Value = $true

# --- Row 15: extr6 (was row 13) ---
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 3).Value = 7
$ws.Cells.Item(15, 4).Value = 11
$ws.Cells.Item(15, 5).Value = $false

# --- Row 16: extr7 (was row 14, now brand new row) ---
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "extr7"
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 7
$ws.Cells.Item(16, 5).Value = $false

# --- Row 17: extr8 (was row 15, now brand new row) ---
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "extr8"
$ws.Cells.Item(17, 3).Value = 8
$ws.Cells.Item(17, 4).Value = 5
$ws.Cells.Item(17, 5).Value = $false
